$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Challenges" cell for row 6 (E6) with the new, longer note.
$ws.Range("E6").Value = "*Excel statistics may not accurately reflect seasonal/holiday commits vs issues due to some holidays being on the edge of months (Christmas/New years)`n*`n*"

# Row 6 needs to grow taller to fit the new text.
$ws.Rows(6).RowHeight = 120

# Update the frozen-pane top-left cell and the active selection in the
# bottom-right pane to reflect the new scroll/selection position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F6").Select()
